$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-number-looking Price values to stay as text (matches source data which
# stores these as text strings, e.g. "575.92"), same as how the existing cells are stored.
$textCells = @("D5", "D6", "D11", "D12", "D13", "D14", "D16", "D19", "D22", "D23", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D41", "D42", "D44", "D47", "D48", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "69.872.24"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "3.574.94"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "575.92"
$ws.Range("E5").Value = "  -3.01%  "
$ws.Range("D6").Value = "187.29"
$ws.Range("E6").Value = "  -3.19%  "
$ws.Range("D7").Value = "3.569.12"
$ws.Range("E7").Value = "  -1.92%  "
$ws.Range("E8").Value = "  -3.52%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").Value = "0.649"
$ws.Range("E11").Value = "  -3.75%  "
$ws.Range("D12").Value = "54.91"
$ws.Range("E12").Value = "  -5.40%  "
$ws.Range("D13").Value = "0.0000302"
$ws.Range("E13").Value = "  +3.30%  "
$ws.Range("D14").Value = "9.56"
$ws.Range("E14").Value = "  -3.65%  "
$ws.Range("D15").Value = "4.150.77"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").Value = "19.65"
$ws.Range("E16").Value = "  -1.94%  "
$ws.Range("D17").Value = "3.577.56"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "69.887.53"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").Value = "12.56"
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("E21").Value = "  -3.11%  "
$ws.Range("D22").Value = "488.48"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").Value = "19.37"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("E24").Value = "  -7.55%  "
$ws.Range("D25").Value = "4.41"
$ws.Range("E25").Value = "  -2.65%  "
$ws.Range("D26").Value = "95.36"
$ws.Range("E26").Value = "  +4.72%  "
$ws.Range("D27").Value = "11.33"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("E28").Value = "  -6.31%  "
$ws.Range("D29").Value = "9.30"
$ws.Range("E29").Value = "  -2.97%  "
$ws.Range("D30").Value = "31.73"
$ws.Range("E30").Value = "  -2.96%  "
$ws.Range("D31").Value = "7.60"
$ws.Range("E31").Value = "  -2.71%  "
$ws.Range("D32").Value = "66.67"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").Value = "12.04"
$ws.Range("E33").Value = "  -1.54%  "
$ws.Range("E34").Value = "  -5.16%  "
$ws.Range("D35").Value = "568.86"
$ws.Range("E35").Value = "  -9.13%  "
$ws.Range("D36").Value = "3.20"
$ws.Range("E36").Value = "  +13.16%  "
$ws.Range("D37").Value = "38.75"
$ws.Range("E37").Value = "  -4.53%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("E40").Value = "  -3.51%  "
$ws.Range("D41").Value = "3.25"
$ws.Range("E41").Value = "  +8.91%  "
$ws.Range("D42").Value = "3.50"
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("E43").Value = "  -7.89%  "
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").Value = "3.02"
$ws.Range("E44").Value = "  -3.86%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "3.256.33"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "9.66"
$ws.Range("E47").Value = "  +4.04%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "3.41"
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("E49").Value = "  -1.90%  "
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "3.19"
$ws.Range("E51").Value = "  -3.33%  "
